{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph right before it) that used to follow the\n// \"LOB1011: Eletricidade Aplicada (Requisito fraco)\" requirements line.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the anchor paragraph (\"LOB1011: ...\") so we only touch the\n// footer paragraphs that immediately follow it.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.trim() === \"LOB1011: Eletricidade Aplicada (Requisito fraco)\") {\n    anchorIndex = i;\n    break;\n  }\n}\n\nconst blankPara = anchorIndex !== -1 ? items[anchorIndex + 1] : null;\nconst jupiterPara = anchorIndex !== -1 ? items[anchorIndex + 2] : null;\nconst copyrightPara = anchorIndex !== -1 ? items[anchorIndex + 3] : null;\n\nconst matches =\n  anchorIndex !== -1 &&\n  blankPara && blankPara.text.trim() === \"\" &&\n  jupiterPara && jupiterPara.text.trim() === \"Ver no Jupiter Salvar em pdf Salvar em docx\" &&\n  copyrightPara && copyrightPara.text.trim().startsWith(\"\u00a9 2020\");\n\nif (matches) {\n  // Delete from the end backwards so earlier indices stay valid.\n  copyrightPara.delete();\n  jupiterPara.delete();\n  blankPara.delete();\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph right before it) that used to follow the\n# \"LOB1011: Eletricidade Aplicada (Requisito fraco)\" requirements line.\n\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\n\n# Locate the anchor paragraph (\"LOB1011: ...\") so we only touch the\n# footer paragraphs that immediately follow it.\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq \"LOB1011: Eletricidade Aplicada (Requisito fraco)\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ne -1) {\n    $blankIndex = $anchorIndex + 1\n    $jupiterIndex = $anchorIndex + 2\n    $copyrightIndex = $anchorIndex + 3\n\n    if (($copyrightIndex + 1) -le $count) {\n        $blankText = $d.Paragraphs.Item($blankIndex).Range.Text.Trim()\n        $jupiterText = $d.Paragraphs.Item($jupiterIndex).Range.Text.Trim()\n        $copyrightText = $d.Paragraphs.Item($copyrightIndex).Range.Text.Trim()\n    } else {\n        $blankText = $null\n        $jupiterText = $null\n        $copyrightText = \"\"\n    }\n\n    if ($blankText -eq \"\" -and\n        $jupiterText -eq \"Ver no Jupiter Salvar em pdf Salvar em docx\" -and\n        $copyrightText.StartsWith(\"\u00a9 2020\")) {\n\n        $startPara = $d.Paragraphs.Item($blankIndex)\n        # One past the last paragraph to delete, so its paragraph mark is\n        # preserved and the following paragraphs are left untouched.\n        $afterPara = $d.Paragraphs.Item($copyrightIndex + 1)\n\n        $r = $d.Range($startPara.Range.Start, $afterPara.Range.Start)\n        $r.Delete()\n    }\n}\n"}
